$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-less approach: Excel auto-converts a purely numeric-looking string
# typed into a cell into a Number. The source cells here are plain TEXT
# (inline strings), e.g. "4.73", "18.04 (15.25%)", so to faithfully
# reproduce a text value like "4.89" we briefly format the cell as Text,
# assign the value, then restore the default "Normal" style so no stray
# number-format style sticks around on the cell.

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Rows 19-24: PLK Regular price (I) -> "4.89"; PLK Percentage Tiered Prices (J) -> "5:3,10:6,20:9,30:12,40:15"
foreach ($r in 19..24) {
    Set-TextValue $ws.Cells.Item($r, 9) "4.89"
    $ws.Cells.Item($r, 10).Value = "5:3,10:6,20:9,30:12,40:15"
}

# Rows 25-29: PLK Regular price (I) -> "21.29"; PLK Percentage Tiered Prices (J) -> "6:3,12:6,24:9,26:12,48:15"
foreach ($r in 25..29) {
    Set-TextValue $ws.Cells.Item($r, 9) "21.29"
    $ws.Cells.Item($r, 10).Value = "6:3,12:6,24:9,26:12,48:15"
}

# Rows 30-34: PLK Regular price (I) -> "21.49"; PLK Percentage Tiered Prices (J) -> "6:3,12:6,24:9,26:12,40:15"
foreach ($r in 30..34) {
    Set-TextValue $ws.Cells.Item($r, 9) "21.49"
    $ws.Cells.Item($r, 10).Value = "6:3,12:6,24:9,26:12,40:15"
}

# Row 35: PLK Sale price (H) -> "13.49"
Set-TextValue $ws.Cells.Item(35, 8) "13.49"
